# edit.ps1 - applies the "updated README and binder" content changes
# to "Problem Formulation & EDA.pptx" per the target diff.
#
# Slide-number note: the pc:sldMk cId/sldId values in the diff refer to
# the deck's internal slide IDs; mapped to 1-based positions in the
# actual slide show these are:
#   sldId 273 -> slide position 7  (shape text tweaks)
#   sldId 277 -> slide position 13 (new "interact" bullet)
#   sldId 264 -> slide position 17 (widen the process-arrow caption box)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 7 (sldId 273): "The 1-liner statement"
#   - extend the price-adjustments sentence
#   - extend the sales-performance bullet
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(3)
$tr7 = $shp7.TextFrame.TextRange

$para4 = $tr7.Paragraphs(4, 1)
$run4b = $para4.Runs(2, 1)
$run4b.Text = "make informed decisions on price adjustments based on risk profile"

$para6 = $tr7.Paragraphs(6, 1)
$run6a = $para6.Runs(1, 1)
$run6a.Text = "The sales performance and profitability"

# ---------------------------------------------------------------------
# Slide 13 (sldId 277): "Framework for EDA"
#   - add a new "interact" bullet under Toolset, right after
#     "Scatter matrix" and before "Correlation (Pearson, Kendall, Spearman)"
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item(2)
$tr13 = $shp13.TextFrame.TextRange

$scatterPara = $tr13.Paragraphs(5, 1)
$scatterPara.InsertAfter("`rinteract")

# ---------------------------------------------------------------------
# Slide 17 (sldId 264): "Summary"
#   - widen the "Define problem -> collect -> ..." caption textbox
# ---------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$shp17 = $s17.Shapes.Item(2)
$shp17.Width = 887
